$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3107.4075

$ws.Range("H79").Value = 3107.4075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 38471644
$ws.Range("I45").Value = 45465932
$ws.Range("J45").Value = 3057
$ws.Range("K45").Value = 45465932
$ws.Range("L45").Value = 3057
$ws.Range("M45").Value = -45465555
$ws.Range("N45").Value = -3811

$ws.Range("H61").Value = 2711.348
$ws.Range("I61").Value = 2335.0625
$ws.Range("J61").Value = 3571.4285
$ws.Range("K61").Value = 2335.0625
$ws.Range("L61").Value = 3571.4285
$ws.Range("M61").Value = -2123.0625
$ws.Range("N61").Value = -3995.4285

$ws.Range("H63").Value = 2468.254
$ws.Range("J63").Value = 2783.3333
$ws.Range("L63").Value = 2783.3333
$ws.Range("N63").Value = -4155.3333

$ws.Range("H66").Value = 2468.254
$ws.Range("J66").Value = 2783.3333
$ws.Range("L66").Value = 13916.6665
$ws.Range("N66").Value = -20780.6665

$ws.Range("H76").Value = 24166.666
$ws.Range("J76").Value = 24166.666
$ws.Range("L76").Value = 24166.666
$ws.Range("N76").Value = -24842.666

$ws.Range("H79").Value = 24166.666
$ws.Range("J79").Value = 24166.666
$ws.Range("L79").Value = 24166.666
$ws.Range("N79").Value = -26506.666

$ws.Range("H122").Value = 4631017
$ws.Range("I122").Value = 5953684
$ws.Range("J122").Value = 1683.3334
$ws.Range("K122").Value = 17861052
$ws.Range("L122").Value = 5050.0002
$ws.Range("M122").Value = -17858602
$ws.Range("N122").Value = -9950.0002

$ws.Range("H136").Value = 2711.348
$ws.Range("I136").Value = 2335.0625
$ws.Range("J136").Value = 3571.4285
$ws.Range("K136").Value = 7005.1875
$ws.Range("L136").Value = 10714.2855
$ws.Range("M136").Value = -4455.1875
$ws.Range("N136").Value = -15814.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1517281.5
$ws.Range("I105").Value = 2274592.2
$ws.Range("J105").Value = 2660
$ws.Range("K105").Value = 2274592.2
$ws.Range("L105").Value = 2660
$ws.Range("M105").Value = -2272845.2
$ws.Range("N105").Value = -6154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14495310
$ws.Range("I31").Value = 22728250
$ws.Range("J31").Value = 5333.28
$ws.Range("K31").Value = 22728250
$ws.Range("L31").Value = 5333.28
$ws.Range("M31").Value = -22727955
$ws.Range("N31").Value = -5923.28

$ws.Range("H34").Value = 14495310
$ws.Range("I34").Value = 22728250
$ws.Range("J34").Value = 5333.28
$ws.Range("K34").Value = 22728250
$ws.Range("L34").Value = 5333.28
$ws.Range("M34").Value = -22728048
$ws.Range("N34").Value = -5737.28

$ws.Range("H58").Value = 1691.9138
$ws.Range("I58").Value = 1270.4062
$ws.Range("J58").Value = 2210.6924
$ws.Range("K58").Value = 1270.4062
$ws.Range("L58").Value = 2210.6924
$ws.Range("M58").Value = -1067.4062
$ws.Range("N58").Value = -2616.6924

$ws.Range("H62").Value = 3108
$ws.Range("I62").Value = 2820
$ws.Range("J62").Value = 3348
$ws.Range("K62").Value = 2820
$ws.Range("L62").Value = 3348
$ws.Range("M62").Value = -2196
$ws.Range("N62").Value = -4596

$ws.Range("H65").Value = 3108
$ws.Range("I65").Value = 2820
$ws.Range("J65").Value = 3348
$ws.Range("K65").Value = 14100
$ws.Range("L65").Value = 16740
$ws.Range("M65").Value = -10980
$ws.Range("N65").Value = -22980

$ws.Range("H136").Value = 1691.9138
$ws.Range("I136").Value = 1270.4062
$ws.Range("J136").Value = 2210.6924
$ws.Range("K136").Value = 3811.2186
$ws.Range("L136").Value = 6632.0772
$ws.Range("M136").Value = -1261.2186
$ws.Range("N136").Value = -11732.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 771.82355
$ws.Range("J5").Value = 1170.8572
$ws.Range("L5").Value = 3512.5716
$ws.Range("N5").Value = -3736.5716

$ws.Range("H122").Value = 900.5833
$ws.Range("I122").Value = 729.2222
$ws.Range("J122").Value = 1071.9445
$ws.Range("K122").Value = 6562.999800000001
$ws.Range("L122").Value = 9647.5005
$ws.Range("M122").Value = -4112.999800000001
$ws.Range("N122").Value = -14547.5005

$ws.Range("H131").Value = 9654846
$ws.Range("J131").Value = 49985.242
$ws.Range("L131").Value = 149955.726
$ws.Range("N131").Value = -160035.726

$ws.Range("H135").Value = 771.82355
$ws.Range("J135").Value = 1170.8572
$ws.Range("L135").Value = 10537.7148
$ws.Range("N135").Value = -15607.7148

$ws.Range("H136").Value = 2907994.2
$ws.Range("I136").Value = 62502804
$ws.Range("J136").Value = 930.4878
$ws.Range("K136").Value = 187508412
$ws.Range("L136").Value = 2791.4634
$ws.Range("M136").Value = -187503312
$ws.Range("N136").Value = -12991.4634

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5127.7617
$ws.Range("I70").Value = 4981.3887
$ws.Range("J70").Value = 6006
$ws.Range("K70").Value = 4981.3887
$ws.Range("L70").Value = 6006
$ws.Range("M70").Value = -4711.3887
$ws.Range("N70").Value = -6546

$ws.Range("H73").Value = 5127.7617
$ws.Range("I73").Value = 4981.3887
$ws.Range("J73").Value = 6006
$ws.Range("K73").Value = 4981.3887
$ws.Range("L73").Value = 6006
$ws.Range("M73").Value = -4045.3887
$ws.Range("N73").Value = -7878

$ws.Range("H80").Value = 56288.05
$ws.Range("I80").Value = 2601
$ws.Range("J80").Value = 75462
$ws.Range("K80").Value = 2601
$ws.Range("L80").Value = 75462
$ws.Range("M80").Value = -1603
$ws.Range("N80").Value = -77458

$ws.Range("H83").Value = 56288.05
$ws.Range("I83").Value = 2601
$ws.Range("J83").Value = 75462
$ws.Range("K83").Value = 13005
$ws.Range("L83").Value = 377310
$ws.Range("M83").Value = -8013
$ws.Range("N83").Value = -387294

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6581.1816
$ws.Range("I93").Value = 7446.2354
$ws.Range("J93").Value = 3640
$ws.Range("K93").Value = 7446.2354
$ws.Range("L93").Value = 3640
$ws.Range("M93").Value = -6198.2354
$ws.Range("N93").Value = -6136

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1326.4509
$ws.Range("I132").Value = 1131.4103
$ws.Range("J132").Value = 1960.3334
$ws.Range("K132").Value = 3394.2309
$ws.Range("L132").Value = 5881.0002
$ws.Range("M132").Value = -864.2309
$ws.Range("N132").Value = -10941.0002
